# Regenerate merged AHB files:
#  1. Rename the header labels from the "_old"/"_new" suffix convention
#     used by the previous merge (FV2404 vs FV2410) to the explicit
#     version-tagged names.
#  2. Turn the populated range into a real Excel Table (ListObject) so the
#     header row carries filters and the table is addressable by name.
#  3. Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -----------------------------------------
$usedCols = $ws.UsedRange.Columns.Count()
for ($c = 1; $c -le $usedCols; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $label = $cell.Value()
    if ($label -ne $null) {
        $newLabel = $label.Replace("_old", "_FV2404").Replace("_new", "_FV2410")
        if ($newLabel -ne $label) {
            $cell.Value = $newLabel
        }
    }
}

# --- 2. Convert the range into a Table (ListObject) -----------------------
$usedRows = $ws.UsedRange.Rows.Count()
$firstCell = $ws.Cells.Item(1, 1)
$lastCell = $ws.Cells.Item($usedRows, $usedCols)
$tableRange = $ws.Range($firstCell, $lastCell)
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1, "")
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the top (header) row ---------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
